$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0000006788371158939305
$ws.Range("E2").Value = 0.0000006788371158939305

$ws.Range("D3").Value = 0.00002568514543710756
$ws.Range("E3").Value = 0.00002568514543710756

$ws.Range("D4").Value = 0.000000000000000001049925520631145
$ws.Range("E4").Value = 0.000000000000000001049925520631145

$ws.Range("D5").Value = 0.000000000000000000110775702984492
$ws.Range("E5").Value = 0.000000000000000000110775702984492

$ws.Range("D6").Value = 0.0007228411622045652
$ws.Range("E6").Value = 0.0007228411622045652

$ws.Range("D7").Value = 0.000003894527830186534
$ws.Range("E7").Value = 0.9999961054721698

$ws.Range("D8").Value = 0.7801052370727506
$ws.Range("E8").Value = 0.2198947629272494

$ws.Range("D9").Value = 0.9999998208745907
$ws.Range("E9").Value = 0.0000001791254092831096

$ws.Range("D10").Value = 0.9999782060605386
$ws.Range("E10").Value = 0.00002179393946144526

$ws.Range("D11").Value = 0.9993007604893783
$ws.Range("E11").Value = 0.0006992395106216742
$ws.Range("F11").Value = 1.270573735237122
